# Generate Report for Handback
# Update the timestamp values on the Overview, zh-cn, and de-de sheets
# to reflect the new handback/handoff generation times.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row
$wsOverview.Range("G2").Value = "2016-09-04 05:11:12"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn.Range("H2").Value = "2016-09-04 05:11:07"
$wsZhCn.Range("K2").Value = "2016-09-04 05:11:31"

# de-de sheet: Correspond Handback DateTime (K2)
$wsDeDe.Range("K2").Value = "2016-09-04 05:11:39"
